$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'313.92"
$ws.Range("E2").Value = "'3.28%"
$ws.Range("D3").Value = "'34.99"
$ws.Range("E3").Value = "'-2.23%"
$ws.Range("D4").Value = "'5.120"
$ws.Range("E4").Value = "'1.10%"
$ws.Range("D5").Value = "'0.08148"
$ws.Range("E5").Value = "'3.47%"
$ws.Range("D6").Value = "'2.126"
$ws.Range("E6").Value = "'0.72%"
$ws.Range("D7").Value = "'4.147"
$ws.Range("E7").Value = "'0.40%"
$ws.Range("D8").Value = "'7.950"
$ws.Range("E8").Value = "'0.08%"
$ws.Range("D9").Value = "'0.9328"
$ws.Range("E9").Value = "'1.26%"
$ws.Range("D10").Value = "'0.1033"
$ws.Range("E10").Value = "'6.28%"
$ws.Range("D11").Value = "'0.1929"
$ws.Range("E11").Value = "'4.73%"
$ws.Range("D12").Value = "'0.09020"
$ws.Range("E12").Value = "'4.79%"
$ws.Range("D13").Value = "'0.03656"
$ws.Range("E13").Value = "'2.12%"
$ws.Range("D14").Value = "'0.09888"
$ws.Range("E14").Value = "'-0.61%"
$ws.Range("D15").Value = "'0.001436"
$ws.Range("E15").Value = "'-0.02%"
$ws.Range("D16").Value = "'0.005852"
$ws.Range("E16").Value = "'2.04%"
$ws.Range("D17").Value = "'3.468"
$ws.Range("E17").Value = "'0.19%"
$ws.Range("E18").Value = "'5.33%"
$ws.Range("D19").Value = "'0.3409"
$ws.Range("E19").Value = "'0.85%"
$ws.Range("D20").Value = "'0.1331"
$ws.Range("E20").Value = "'-1.32%"
$ws.Range("D21").Value = "'5.116"
$ws.Range("E21").Value = "'-1.29%"
$ws.Range("E22").Value = "'0.24%"
$ws.Range("D23").Value = "'0.04552"
$ws.Range("E23").Value = "'0.06%"
$ws.Range("D24").Value = "'0.001249"
$ws.Range("E24").Value = "'0.94%"
$ws.Range("D25").Value = "'0.004697"
$ws.Range("E25").Value = "'-3.58%"
$ws.Range("D26").Value = "'0.0001253"
$ws.Range("E26").Value = "'-3.98%"
$ws.Range("D27").Value = "'0.0004508"
$ws.Range("E27").Value = "'-5.33%"
$ws.Range("D39").Value = "'0.01948"
$ws.Range("E39").Value = "'5.40%"
$ws.Range("D40").Value = "'0.04899"
$ws.Range("E40").Value = "'4.33%"
$ws.Range("D41").Value = "'0.007618"
$ws.Range("E41").Value = "'-3.69%"
$ws.Range("D42").Value = "'0.1388"
$ws.Range("E42").Value = "'-0.34%"
$ws.Range("D43").Value = "'0.007878"
$ws.Range("E43").Value = "'3.74%"
$ws.Range("D44").Value = "'0.002100"
$ws.Range("E44").Value = "'-4.45%"
$ws.Range("D45").Value = "'0.01180"
$ws.Range("E45").Value = "'4.88%"
$ws.Range("D46").Value = "'0.00006756"
$ws.Range("E46").Value = "'7.62%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("D48").Value = "'186.61"
$ws.Range("E48").Value = "'269.24%"
$ws.Range("D49").Value = "'0.001703"
$ws.Range("E49").Value = "'-10.60%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'-0.20%"
